# Apply report_v2 -> report_v3 changes.
#
# The edits are small, surgical OOXML-level tweaks (removing stray
# <w:rFonts w:hint="eastAsia"/> paragraph-mark runs left over from a
# Korean IME, collapsing a bookmark-split "TeamViewer" back into one
# run, trimming three blank trailing paragraphs down to a bookmark, and
# marking a style as semiHidden). Word's object model doesn't expose a
# "delete this exact pPr/rPr" verb, so we do it the way a macro author
# would via Find & Replace on the raw WordOpenXML package text, which
# round-trips losslessly through this document.
#
# NOTE: this runtime's function/param binding mangles large strings
# passed via named parameters, so the replacements are inlined here
# rather than wrapped in a helper function.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# 1) RDP section: drop the stray paragraph-mark <w:rPr><w:rFonts hint=eastAsia/></w:rPr>
#    on the (otherwise empty) paragraph that holds the screenshot.
$old1 = '<w:p w14:paraId="33265938" w14:textId="339EB3A2" w:rsidR="00337340" w:rsidRDefault="00337340" w:rsidP="00992853"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r>'
$new1 = '<w:p w14:paraId="33265938" w14:textId="339EB3A2" w:rsidR="00337340" w:rsidRDefault="00337340" w:rsidP="00992853"><w:r>'
if ($xml.IndexOf($old1) -lt 0) { throw "pattern 1 (RDP screenshot paragraph) not found" }
$xml = $xml.Replace($old1, $new1)

# 2) Trailing blank paragraph in a table cell: was pPr/rPr hint=eastAsia only, becomes a bare <w:p/>.
$old2 = '<w:p w14:paraId="48762548" w14:textId="6EFACE87" w:rsidR="0099526F" w:rsidRDefault="0099526F" w:rsidP="001D6A7B"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>'
$new2 = '<w:p w14:paraId="48762548" w14:textId="6EFACE87" w:rsidR="0099526F" w:rsidRDefault="0099526F" w:rsidP="001D6A7B"/>'
if ($xml.IndexOf($old2) -lt 0) { throw "pattern 2 (blank trailing cell paragraph) not found" }
$xml = $xml.Replace($old2, $new2)

# 3) VNC section: same stray paragraph-mark rFonts fix as (1), on the paragraph before that screenshot.
$old3 = '<w:p w14:paraId="75A98642" w14:textId="5173DE67" w:rsidR="007A74D3" w:rsidRPr="00AA6028" w:rsidRDefault="00284ACB" w:rsidP="0024495F"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r>'
$new3 = '<w:p w14:paraId="75A98642" w14:textId="5173DE67" w:rsidR="007A74D3" w:rsidRPr="00AA6028" w:rsidRDefault="00284ACB" w:rsidP="0024495F"><w:r>'
if ($xml.IndexOf($old3) -lt 0) { throw "pattern 3 (VNC screenshot paragraph) not found" }
$xml = $xml.Replace($old3, $new3)

# 4) "Tea" + _GoBack bookmark + "mViewer" (an artifact of mid-word cursor tracking) -> one "TeamViewer" run.
$old4 = '<w:r w:rsidR="00924A5B"><w:t>Tea</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00924A5B"><w:t>mViewer</w:t></w:r>'
$new4 = '<w:r w:rsidR="00924A5B"><w:t>TeamViewer</w:t></w:r>'
if ($xml.IndexOf($old4) -lt 0) { throw "pattern 4 (TeamViewer bookmark split) not found" }
$xml = $xml.Replace($old4, $new4)

# 5) Heading paragraph mark: drop the stray <w:rFonts hint=eastAsia/> from the b/bCs rPr.
$old5 = '<w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr></w:pPr>'
$new5 = '<w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>'
if ($xml.IndexOf($old5) -lt 0) { throw "pattern 5 (numPr heading rFonts) not found" }
$xml = $xml.Replace($old5, $new5)

# 6) The RDP-vs-VNC comparison paragraph gains a paragraph-mark <w:rFonts hint=eastAsia/>.
$old6 = 'w:rsidRDefault="004E7992" w:rsidP="004F3FD9"><w:pPr><w:ind w:firstLineChars="100" w:firstLine="200"/></w:pPr>'
$new6 = 'w:rsidRDefault="004E7992" w:rsidP="004F3FD9"><w:pPr><w:ind w:firstLineChars="100" w:firstLine="200"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>'
if ($xml.IndexOf($old6) -lt 0) { throw "pattern 6 (firstLineChars paragraph rFonts) not found" }
$xml = $xml.Replace($old6, $new6)

# 7) Collapse the three trailing blank paragraphs (each just a lone space, hint=eastAsia) in the
#    last table cell down to a bookmark placed at the end of the preceding paragraph.
$old7 = '<w:r w:rsidR="00AD39CD"><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="23EBE81E" w14:textId="3F02368C" w:rsidR="00255E1A" w:rsidRDefault="00753CBE" w:rsidP="0008358F"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="0DDDFA67" w14:textId="17DBE77C" w:rsidR="00753CBE" w:rsidRDefault="00753CBE" w:rsidP="0008358F"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="5F9CD034" w14:textId="418B9BFC" w:rsidR="00255E1A" w:rsidRPr="00CA1D40" w:rsidRDefault="00255E1A" w:rsidP="0008358F"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>'
$new7 = '<w:r w:rsidR="00AD39CD"><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
if ($xml.IndexOf($old7) -lt 0) { throw "pattern 7 (trailing blank paragraphs collapse) not found" }
$xml = $xml.Replace($old7, $new7)

# 8) Mark the "Default Paragraph Font" character style as semiHidden.
$old8 = '<w:style w:type="character" w:default="1" w:styleId="a0"><w:name w:val="Default Paragraph Font"/><w:uiPriority w:val="1"/><w:unhideWhenUsed/></w:style>'
$new8 = '<w:style w:type="character" w:default="1" w:styleId="a0"><w:name w:val="Default Paragraph Font"/><w:uiPriority w:val="1"/><w:semiHidden/><w:unhideWhenUsed/></w:style>'
if ($xml.IndexOf($old8) -lt 0) { throw "pattern 8 (Default Paragraph Font semiHidden) not found" }
$xml = $xml.Replace($old8, $new8)

$d.WordOpenXML = $xml

Write-Host "Applied 8 replacements successfully."
